$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 230.30435
$ws.Range("I33").Value = 230.30435
$ws.Range("K33").Value = 230.30435
$ws.Range("M33").Value = -1.304349999999999
$ws.Range("H74").Value = 3659.625
$ws.Range("I74").Value = 2759
$ws.Range("K74").Value = 2759
$ws.Range("M74").Value = -1823
$ws.Range("H77").Value = 3659.625
$ws.Range("I77").Value = 2759
$ws.Range("K77").Value = 13795
$ws.Range("M77").Value = -9115
$ws.Range("H106").Value = 1486.875
$ws.Range("I106").Value = 1048.75
$ws.Range("K106").Value = 1048.75
$ws.Range("M106").Value = -417.75
$ws.Range("H139").Value = 72423.75
$ws.Range("J139").Value = 72423.75
$ws.Range("L139").Value = 72423.75
$ws.Range("N139").Value = -82703.75
$ws.Range("H140").Value = 96079.234
$ws.Range("J140").Value = 96290
$ws.Range("L140").Value = 96290
$ws.Range("N140").Value = -106650

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 3000
$ws.Range("I6").Value = 3000
$ws.Range("K6").Value = 3000
$ws.Range("M6").Value = -2827
$ws.Range("H63").Value = 13425
$ws.Range("J63").Value = 9900
$ws.Range("L63").Value = 9900
$ws.Range("N63").Value = -11272
$ws.Range("H66").Value = 13425
$ws.Range("J66").Value = 9900
$ws.Range("L66").Value = 49500
$ws.Range("N66").Value = -56364
$ws.Range("H102").Value = 1251.9286
$ws.Range("I102").Value = 1117.4615
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 1117.4615
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = 504.5385000000001
$ws.Range("N102").Value = -6244
$ws.Range("H132").Value = 871964.9
$ws.Range("I132").Value = 1430471.4
$ws.Range("J132").Value = 3176.889
$ws.Range("K132").Value = 4291414.199999999
$ws.Range("L132").Value = 9530.667000000001
$ws.Range("M132").Value = -4288884.199999999
$ws.Range("N132").Value = -14590.667
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H138").Value = 61313.332
$ws.Range("J138").Value = 61313.332
$ws.Range("L138").Value = 61313.332
$ws.Range("N138").Value = -71593.33199999999
$ws.Range("H139").Value = 84707.88
$ws.Range("J139").Value = 84707.88
$ws.Range("L139").Value = 84707.88
$ws.Range("N139").Value = -94987.88
$ws.Range("H141").Value = 57810
$ws.Range("J141").Value = 57810
$ws.Range("L141").Value = 57810
$ws.Range("N141").Value = -68170

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 10800.444
$ws.Range("I8").Value = 12075.5
$ws.Range("J8").Value = 600
$ws.Range("K8").Value = 12075.5
$ws.Range("L8").Value = 600
$ws.Range("M8").Value = -11935.5
$ws.Range("N8").Value = -880
$ws.Range("H74").Value = 91780
$ws.Range("J74").Value = 91780
$ws.Range("L74").Value = 91780
$ws.Range("N74").Value = -93652
$ws.Range("H77").Value = 91780
$ws.Range("J77").Value = 91780
$ws.Range("L77").Value = 275340
$ws.Range("N77").Value = -284700
$ws.Range("H132").Value = 76180.836
$ws.Range("J132").Value = 76180.836
$ws.Range("L132").Value = 76180.836
$ws.Range("N132").Value = -86300.836
$ws.Range("H135").Value = 76346.664
$ws.Range("J135").Value = 76346.664
$ws.Range("L135").Value = 76346.664
$ws.Range("N135").Value = -86486.664
$ws.Range("H138").Value = 49446.875
$ws.Range("J138").Value = 49446.875
$ws.Range("L138").Value = 49446.875
$ws.Range("N138").Value = -59726.875

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H137").Value = 44926.668
$ws.Range("J137").Value = 74780
$ws.Range("L137").Value = 74780
$ws.Range("N137").Value = -84980
$ws.Range("H138").Value = 58797.855
$ws.Range("J138").Value = 58797.855
$ws.Range("L138").Value = 58797.855
$ws.Range("N138").Value = -69077.85500000001
$ws.Range("H140").Value = 72980
$ws.Range("J140").Value = 72980
$ws.Range("L140").Value = 72980
$ws.Range("N140").Value = -83340

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 4231.567
$ws.Range("I113").Value = 594.3333
$ws.Range("K113").Value = 1782.9999
$ws.Range("M113").Value = 387.0001
$ws.Range("H129").Value = 1853696.2
$ws.Range("J129").Value = 2085312.4
$ws.Range("L129").Value = 6255937.199999999
$ws.Range("N129").Value = -6265937.199999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 8003267
$ws.Range("I3").Value = 834999.9399999999
$ws.Range("K3").Value = 834999.9399999999
$ws.Range("M3").Value = -834883.9399999999
$ws.Range("H17").Value = 2269.6667
$ws.Range("I17").Value = 2000
$ws.Range("J17").Value = 2404.5
$ws.Range("K17").Value = 2000
$ws.Range("L17").Value = 2404.5
$ws.Range("M17").Value = -1832
$ws.Range("N17").Value = -2740.5
$ws.Range("H41").Value = 4017
$ws.Range("J41").Value = 10000
$ws.Range("L41").Value = 10000
$ws.Range("N41").Value = -10710
$ws.Range("H113").Value = 2160
$ws.Range("J113").Value = 2800
$ws.Range("L113").Value = 2800
$ws.Range("N113").Value = -7140
$ws.Range("H132").Value = 3140
$ws.Range("I132").Value = 2065.4285
$ws.Range("J132").Value = 4080.25
$ws.Range("K132").Value = 6196.2855
$ws.Range("L132").Value = 12240.75
$ws.Range("M132").Value = -3666.2855
$ws.Range("N132").Value = -17300.75
$ws.Range("H133").Value = 62410
$ws.Range("J133").Value = 62410
$ws.Range("L133").Value = 62410
$ws.Range("N133").Value = -72530
$ws.Range("H135").Value = 77574.60000000001
$ws.Range("J135").Value = 77574.60000000001
$ws.Range("L135").Value = 77574.60000000001
$ws.Range("N135").Value = -87714.60000000001
$ws.Range("H138").Value = 59770
$ws.Range("J138").Value = 59770
$ws.Range("L138").Value = 59770
$ws.Range("N138").Value = -70050
$ws.Range("H141").Value = 48673.75
$ws.Range("J141").Value = 48673.75
$ws.Range("L141").Value = 48673.75
$ws.Range("N141").Value = -59033.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1499.909
$ws.Range("I68").Value = 1312.375
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 1312.375
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = -563.375
$ws.Range("N68").Value = -3498
$ws.Range("H71").Value = 1499.909
$ws.Range("I71").Value = 1312.375
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 6561.875
$ws.Range("L71").Value = 10000
$ws.Range("M71").Value = -2817.875
$ws.Range("N71").Value = -17488
$ws.Range("H134").Value = 76660
$ws.Range("J134").Value = 76660
$ws.Range("L134").Value = 76660
$ws.Range("N134").Value = -86800
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("H138").Value = 76960
$ws.Range("J138").Value = 76960
$ws.Range("L138").Value = 76960
$ws.Range("N138").Value = -87240
$ws.Range("H139").Value = 61827.273
$ws.Range("J139").Value = 61827.273
$ws.Range("L139").Value = 61827.273
$ws.Range("N139").Value = -72107.273
$ws.Range("H141").Value = 46500
$ws.Range("J141").Value = 46500
$ws.Range("L141").Value = 46500
$ws.Range("N141").Value = -56860

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 65313
$ws.Range("J46").Value = 65313
$ws.Range("L46").Value = 65313
$ws.Range("N46").Value = -65775
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H81").Value = 103588.1
$ws.Range("I81").Value = 202375.8
$ws.Range("J81").Value = 4800.4
$ws.Range("K81").Value = 404751.6
$ws.Range("L81").Value = 9600.799999999999
$ws.Range("M81").Value = -403690.6
$ws.Range("N81").Value = -11722.8
$ws.Range("H84").Value = 103588.1
$ws.Range("I84").Value = 202375.8
$ws.Range("J84").Value = 4800.4
$ws.Range("K84").Value = 2023758
$ws.Range("L84").Value = 48004
$ws.Range("M84").Value = -2018454
$ws.Range("N84").Value = -58612
$ws.Range("H134").Value = 65313
$ws.Range("J134").Value = 65313
$ws.Range("L134").Value = 195939
$ws.Range("N134").Value = -201009
$ws.Range("H138").Value = 58474.75
$ws.Range("J138").Value = 58474.75
$ws.Range("L138").Value = 58474.75
$ws.Range("N138").Value = -68754.75
$ws.Range("H140").Value = 40528.89
$ws.Range("J140").Value = 40528.89
$ws.Range("L140").Value = 40528.89
$ws.Range("N140").Value = -50888.89
$ws.Range("H141").Value = 87935.625
$ws.Range("J141").Value = 87935.625
$ws.Range("L141").Value = 87935.625
$ws.Range("N141").Value = -98295.625
